# ARKCORR-18 Added business process definitions for the on enter queue rule.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the rule CONDITION / ACTION template text (row 17) ---
$ws.Range("C17").Value = '((CaseFilePipelineContext) $model.getPipelineContext()).getEnqueueName().equals("$param")'
$ws.Range("D17").Value = '$model.setBusinessProcessName("$param");'

# --- Update the default-rule row (row 18) description text ---
$ws.Range("C18").Value = "Entering Queue Name"
$ws.Range("D18").Value = "The business process to be executed when entering that queue"

# --- Apply the same formatting used by the existing data rows (B18) to the
#     new rule rows below it, then fill in the per-queue rule data ---
$ws.Range("B18").Copy() | Out-Null
$ws.Range("B19:D23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B19").Value = "Intake queue"
$ws.Range("C19").Value = "Intake"
$ws.Range("D19").Value = "correspondence-extension-intake-process"

$ws.Range("B20").Value = "Fulfill queue"
$ws.Range("C20").Value = "Fulfill"
$ws.Range("D20").Value = "correspondence-extension-fulfill-process"

$ws.Range("B21").Value = "Supervisor Approval queue"
$ws.Range("C21").Value = "Supervisor Approval"
$ws.Range("D21").Value = "correspondence-extension-supervisor-approval-process"

$ws.Range("B22").Value = "Executive Approval queue"
$ws.Range("C22").Value = "Executive Approval"
$ws.Range("D22").Value = "correspondence-extension-executive-approval-process"

$ws.Range("B23").Value = "Release queue"
$ws.Range("C23").Value = "Release"
$ws.Range("D23").Value = "correspondence-extension-release-process"

# --- Refresh the used range / selection so the sheet view matches the new extent ---
$ws.Range("D23").Select() | Out-Null

Write-Host "Edit complete"
